# Insert a new weekly price record as row 8, pushing the existing
# rows 8-48 down to rows 9-49 (dimension grows from A1:R48 to A1:R49).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 8..48 down by one row.
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row 8 with the new record.
$ws.Cells.Item(8, 1).Value  = 1
$ws.Cells.Item(8, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(8, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(8, 4).Value  = 44473
$ws.Cells.Item(8, 5).Value  = 15
$ws.Cells.Item(8, 6).Value  = 100112038
$ws.Cells.Item(8, 7).Value  = "Cebollín baby"
$ws.Cells.Item(8, 8).Value  = "Sin especificar"
$ws.Cells.Item(8, 9).Value  = "Primera"
$ws.Cells.Item(8, 10).Value = 300
$ws.Cells.Item(8, 11).Value = 950
$ws.Cells.Item(8, 12).Value = 1000
$ws.Cells.Item(8, 13).Value = 975
$ws.Cells.Item(8, 14).Value = "`$/paquete 1,5 a 2 kilos"
$ws.Cells.Item(8, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(8, 16).Value = 488
$ws.Cells.Item(8, 17).Value = 2
$ws.Cells.Item(8, 18).Value = "Hortaliza"
